$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text like "42.933.63" which Excel may otherwise
# try to parse as a number (losing trailing zeros / precision), so force
# the cells to be treated as plain text before writing the values.
$dCells = "D2","D3","D5","D7","D10","D12","D14","D15","D16","D18","D21","D22","D23","D24","D26","D27","D31","D36","D38","D39","D40","D44","D45","D46","D47","D48","D51"
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.933.63"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").Value = "2.218.59"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "257.07"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("D7").Value = "77.33"
$ws.Range("E7").Value = "  +2.90%  "

$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("D10").Value = "42.44"
$ws.Range("E10").Value = "  +3.66%  "

$ws.Range("E11").Value = "  -2.05%  "

$ws.Range("D12").Value = "7.06"
$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("E13").Value = "  +1.18%  "

$ws.Range("D14").Value = "2.551.32"
$ws.Range("E14").Value = "  -1.37%  "

$ws.Range("D15").Value = "14.56"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").Value = "2.222.76"
$ws.Range("E16").Value = "  -0.93%  "

$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").Value = "42.921.67"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("E19").Value = "  -1.40%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  +0.42%  "

$ws.Range("D22").Value = "2.24"
$ws.Range("E22").Value = "  +2.72%  "

$ws.Range("D23").Value = "230.93"
$ws.Range("E23").Value = "  +0.47%  "

$ws.Range("D24").Value = "9.35"
$ws.Range("E24").Value = "  -5.92%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").Value = "43.14"
$ws.Range("E26").Value = "  +11.49%  "

$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  -0.48%  "

$ws.Range("E28").Value = "  -2.31%  "

$ws.Range("E29").Value = "  -1.50%  "

$ws.Range("E30").Value = "  +3.92%  "

$ws.Range("D31").Value = "172.91"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("E33").Value = "  +10.62%  "

$ws.Range("E34").Value = "  -0.73%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "0.0362"
$ws.Range("E36").Value = "  +9.05%  "

$ws.Range("E37").Value = "  -0.20%  "

$ws.Range("D38").Value = "4.36"
$ws.Range("E38").Value = "  +0.33%  "

$ws.Range("D39").Value = "12.96"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +18.52%  "

$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").Value = "60.26"
$ws.Range("E44").Value = "  +1.42%  "

$ws.Range("D45").Value = "103.17"
$ws.Range("E45").Value = "  -1.73%  "

$ws.Range("D46").Value = "8.37"
$ws.Range("E46").Value = "  -3.66%  "

$ws.Range("D47").Value = "0.0978"
$ws.Range("E47").Value = "  -1.25%  "

$ws.Range("D48").Value = "0.466"
$ws.Range("E48").Value = "  -3.12%  "

$ws.Range("E49").Value = "  +1.23%  "

$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("D51").Value = "2.438.71"
$ws.Range("E51").Value = "  -0.86%  "
